$wb = $excel.ActiveWorkbook

$oldVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

# Update the "About" sheet
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = "Version: $newVersion"
$wsAbout.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Integra Underground Mine, Australia, M4354, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# Update the "Boundaries and methane sources" sheet build_version column (S2:S27)
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 27; $row++) {
    $wsData.Cells.Item($row, 19).Value = $newVersion
}
